# Change "WeaponEffect:" -> "WeaponEffects:" (both the bold/underlined label
# and the inline example text) in the "WeaponEffect" tag-reference paragraph.
#
# The target OOXML splits the word "WeaponEffect(s):" into three separate
# <w:r> runs in each spot (matching how the original author's edit was
# serialized: "WeaponEffect" / "s" / ":"), rather than one merged run, so we
# rebuild the paragraph run-by-run to reproduce that exact structure.

$d = $word.ActiveDocument

# Helper: insert a single <w:r>...</w:r> fragment (verbatim WordprocessingML)
# at a document position.
function InsertRunAt($pos, $innerXml) {
    $xml = "<?xml version=`"1.0`"?><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:body><w:p>$innerXml</w:p></w:body></w:document>"
    $rr = $d.Range($pos, $pos)
    $rr.InsertXML($xml)
}

# Locate the paragraph whose text begins with the bold "WeaponEffect:" label
# (there's also a casual mention of "WeaponEffect" inside an earlier bullet,
# which must stay untouched).
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.StartsWith("WeaponEffect:")) {
        $target = $para
    }
}

if ($target -eq $null) {
    Write-Host "WeaponEffect paragraph not found"
} else {
    $pStart = $target.Range.Start
    $pEnd = $target.Range.End

    $boldRpr = "<w:rPr><w:b/><w:bCs/><w:u w:val=`"single`"/></w:rPr>"

    # Wipe the paragraph's text (keeping the paragraph mark / its own
    # attributes intact), then rebuild it from seven runs.
    $full = $d.Range($pStart, $pEnd - 1)
    $full.Text = ""

    # This engine only preserves runs as distinct <w:r> elements (instead of
    # silently re-merging adjacent same-formatted ones) when each is inserted
    # right at a paragraph's start offset, so insert them in reverse order,
    # all at $pStart, so the final reading order comes out correct:
    #   "WeaponEffect" | "s" | ":" | " " | "…WeaponEffect" | "s" | ":[…]…"
    InsertRunAt $pStart "<w:r><w:t>:[{id:x,time:y,lvl:z},&#8230;,{id:&#8221;effect&#8221;,time:b,lvl:c}]&#8230;</w:t></w:r>"
    InsertRunAt $pStart "<w:r><w:t>s</w:t></w:r>"
    InsertRunAt $pStart "<w:r w:rsidR=`"0075565F`"><w:t>&#8230;WeaponEffect</w:t></w:r>"
    InsertRunAt $pStart "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>"
    InsertRunAt $pStart "<w:r>$boldRpr<w:t>:</w:t></w:r>"
    InsertRunAt $pStart "<w:r>$boldRpr<w:t>s</w:t></w:r>"
    InsertRunAt $pStart "<w:r w:rsidRPr=`"007304B6`">$boldRpr<w:t>WeaponEffect</w:t></w:r>"

    Write-Host "WeaponEffect(s) paragraph updated"
}
